# Bugfix for alliance speedup:
# The "rights" sheet lists one row per alliance title (archon, general,
# diplomat, quartermaster, supervisor, elite, member). The "diplomat"
# title/role should not exist, so remove its entire row. Excel shifts
# the rows below it up automatically (and the shared-strings table is
# compacted since "diplomat" is no longer referenced anywhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights")

# Row 4 holds the "diplomat" title entry - delete the whole row, shifting
# everything below up by one (matches Excel's Rows(n).Delete behavior).
$ws.Rows(4).Delete() | Out-Null

# Leave the selection where the author's commit left it.
$ws.Range("C8").Select() | Out-Null
